$wb = $excel.ActiveWorkbook

# Rename sheets (task order id suffixes updated to new timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556261222646"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556283006597"
$wb.Worksheets.Item(3).Name = "RS_TO-16512556283016615"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512556283471844"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651255628425457"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512556260839703.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255626105887.csv"
$ws1.Range("B4").Value = "go_stims-16512556261069922.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556261203969.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16512556279107747.csv"
$ws2.Range("B3").Value = "TB-16512556282794943.csv"
$ws2.Range("B4").Value = "OB-16512556273849878.csv"
$ws2.Range("B5").Value = "OB-16512556277380068.csv"
$ws2.Range("B6").Value = "ZB-match_5-1651255626310833.csv"
$ws2.Range("B7").Value = "TB-16512556278168015.csv"
$ws2.Range("B8").Value = "ZB-match_6-16512556262656684.csv"
$ws2.Range("B9").Value = "OB-16512556272850084.csv"
$ws2.Range("B10").Value = "ZB-match_6-16512556268316195.csv"

# Sheet 3 (RS) -- no cell content changes, only the sheet name (already updated above)

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556283155515.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556283036633.csv"
$ws4.Range("B4").Value = "MM_stims-16512556283306842.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255628316546.csv"
$ws4.Range("B6").Value = "MM_stims-16512556283461854.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255628331681.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651255628394624.csv"
$ws5.Range("B3").Value = "SAT_stims-16512556283530426.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512556284097714.csv"
$ws5.Range("B5").Value = "SAT_stims-16512556283788762.csv"
